# Add a third row: A3 = "Test3", formatted like A1/A2 (same style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A1's formatting onto A3 (reuses the existing style rather than
# minting a new cellXf), then set the new cell's value.
$ws.Range("A1").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3").Value = "Test3"
